$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp text in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 4 de Abril de 2020 a las 21:22"

# Update reordered / revised country rows (values reflect final target state)
# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 302919
$ws.Cells.Item(4, 3).Value = 25758
$ws.Cells.Item(4, 4).Value = 14686
$ws.Cells.Item(4, 5).Value = 279990
$ws.Cells.Item(4, 6).Value = 7983
$ws.Cells.Item(4, 7).Value = 839
$ws.Cells.Item(4, 8).Value = 8243

# Row 7: Alemania
$ws.Cells.Item(7, 1).Value = "Alemania"
$ws.Cells.Item(7, 2).Value = 95637
$ws.Cells.Item(7, 3).Value = 4478
$ws.Cells.Item(7, 4).Value = 26400
$ws.Cells.Item(7, 5).Value = 67842
$ws.Cells.Item(7, 6).Value = 3936
$ws.Cells.Item(7, 7).Value = 120
$ws.Cells.Item(7, 8).Value = 1395

# Row 13: Suiza
$ws.Cells.Item(13, 1).Value = "Suiza"
$ws.Cells.Item(13, 2).Value = 20505
$ws.Cells.Item(13, 3).Value = 899
$ws.Cells.Item(13, 4).Value = 4846
$ws.Cells.Item(13, 5).Value = 14993
$ws.Cells.Item(13, 6).Value = 391
$ws.Cells.Item(13, 7).Value = 75
$ws.Cells.Item(13, 8).Value = 666

# Row 16: Canada
$ws.Cells.Item(16, 1).Value = "Canada"
$ws.Cells.Item(16, 2).Value = 13872
$ws.Cells.Item(16, 3).Value = 1497
$ws.Cells.Item(16, 4).Value = 2565
$ws.Cells.Item(16, 5).Value = 11079
$ws.Cells.Item(16, 6).Value = 120
$ws.Cells.Item(16, 7).Value = 20
$ws.Cells.Item(16, 8).Value = 228

# Row 33: Ecuador
$ws.Cells.Item(33, 1).Value = "Ecuador"
$ws.Cells.Item(33, 2).Value = 3465
$ws.Cells.Item(33, 3).Value = 97
$ws.Cells.Item(33, 4).Value = 100
$ws.Cells.Item(33, 5).Value = 3193
$ws.Cells.Item(33, 6).Value = 100
$ws.Cells.Item(33, 7).Value = 27
$ws.Cells.Item(33, 8).Value = 172

# Row 43: Peru
$ws.Cells.Item(43, 1).Value = "Peru"
$ws.Cells.Item(43, 2).Value = 1746
$ws.Cells.Item(43, 3).Value = 151
$ws.Cells.Item(43, 4).Value = 914
$ws.Cells.Item(43, 5).Value = 759
$ws.Cells.Item(43, 6).Value = 88
$ws.Cells.Item(43, 7).Value = 12
$ws.Cells.Item(43, 8).Value = 73

# Row 56: Ucrania
$ws.Cells.Item(56, 1).Value = "Ucrania"
$ws.Cells.Item(56, 2).Value = 1225
$ws.Cells.Item(56, 3).Value = 153
$ws.Cells.Item(56, 4).Value = 25
$ws.Cells.Item(56, 5).Value = 1168
$ws.Cells.Item(56, 6).Value = 16
$ws.Cells.Item(56, 7).Value = 5
$ws.Cells.Item(56, 8).Value = 32

# Row 57: Singapur
$ws.Cells.Item(57, 1).Value = "Singapur"
$ws.Cells.Item(57, 2).Value = 1189
$ws.Cells.Item(57, 3).Value = 75
$ws.Cells.Item(57, 4).Value = 297
$ws.Cells.Item(57, 5).Value = 886
$ws.Cells.Item(57, 6).Value = 24
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = 6

# Row 58: Croacia
$ws.Cells.Item(58, 1).Value = "Croacia"
$ws.Cells.Item(58, 2).Value = 1126
$ws.Cells.Item(58, 3).Value = 47
$ws.Cells.Item(58, 4).Value = 119
$ws.Cells.Item(58, 5).Value = 995
$ws.Cells.Item(58, 6).Value = 39
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 12

# Row 72: Bosnia y Herzegovina
$ws.Cells.Item(72, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(72, 2).Value = 623
$ws.Cells.Item(72, 3).Value = 44
$ws.Cells.Item(72, 4).Value = 31
$ws.Cells.Item(72, 5).Value = 571
$ws.Cells.Item(72, 6).Value = 4
$ws.Cells.Item(72, 7).Value = 4
$ws.Cells.Item(72, 8).Value = 21

# Row 74: Kazajistan
$ws.Cells.Item(74, 1).Value = "Kazajistan"
$ws.Cells.Item(74, 2).Value = 531
$ws.Cells.Item(74, 3).Value = 67
$ws.Cells.Item(74, 4).Value = 36
$ws.Cells.Item(74, 5).Value = 490
$ws.Cells.Item(74, 6).Value = 6
$ws.Cells.Item(74, 7).Value = 2
$ws.Cells.Item(74, 8).Value = 5

# Row 76: Libano
$ws.Cells.Item(76, 1).Value = "Libano"
$ws.Cells.Item(76, 2).Value = 520
$ws.Cells.Item(76, 3).Value = 12
$ws.Cells.Item(76, 4).Value = 54
$ws.Cells.Item(76, 5).Value = 449
$ws.Cells.Item(76, 6).Value = 27
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 17

# Row 85: Costa Rica
$ws.Cells.Item(85, 1).Value = "Costa Rica"
$ws.Cells.Item(85, 2).Value = 435
$ws.Cells.Item(85, 3).Value = 19
$ws.Cells.Item(85, 4).Value = 13
$ws.Cells.Item(85, 5).Value = 420
$ws.Cells.Item(85, 6).Value = 13
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 2

# Row 86: Republica de Chipre
$ws.Cells.Item(86, 1).Value = "Republica de Chipre"
$ws.Cells.Item(86, 2).Value = 426
$ws.Cells.Item(86, 3).Value = 30
$ws.Cells.Item(86, 4).Value = 33
$ws.Cells.Item(86, 5).Value = 382
$ws.Cells.Item(86, 6).Value = 11
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 11

# Row 123: Ruanda
$ws.Cells.Item(123, 1).Value = "Ruanda"
$ws.Cells.Item(123, 2).Value = 102
$ws.Cells.Item(123, 3).Value = 13
$ws.Cells.Item(123, 4).Value = 0
$ws.Cells.Item(123, 5).Value = 102
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 0

# Row 124: Trinidad yTobago
$ws.Cells.Item(124, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(124, 2).Value = 101
$ws.Cells.Item(124, 3).Value = 3
$ws.Cells.Item(124, 4).Value = 1
$ws.Cells.Item(124, 5).Value = 94
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 6

# Row 125: Gibraltar
$ws.Cells.Item(125, 1).Value = "Gibraltar"
$ws.Cells.Item(125, 2).Value = 98
$ws.Cells.Item(125, 3).Value = 3
$ws.Cells.Item(125, 4).Value = 52
$ws.Cells.Item(125, 5).Value = 46
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 0

# Row 126: Paraguay
$ws.Cells.Item(126, 1).Value = "Paraguay"
$ws.Cells.Item(126, 2).Value = 96
$ws.Cells.Item(126, 3).Value = 4
$ws.Cells.Item(126, 4).Value = 12
$ws.Cells.Item(126, 5).Value = 81
$ws.Cells.Item(126, 6).Value = 2
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 3

# Row 143: Polinesia Francesa
$ws.Cells.Item(143, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(143, 2).Value = 40
$ws.Cells.Item(143, 3).Value = 1
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 40
$ws.Cells.Item(143, 6).Value = 1
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 0

# Row 144: Puerto Rico
$ws.Cells.Item(144, 1).Value = "Puerto Rico"
$ws.Cells.Item(144, 2).Value = 39
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 1
$ws.Cells.Item(144, 5).Value = 36
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 2

# Row 145: Zambia
$ws.Cells.Item(145, 1).Value = "Zambia"
$ws.Cells.Item(145, 2).Value = 39
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 2
$ws.Cells.Item(145, 5).Value = 36
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 1
